# [github-321] Fix issue with rounding in DataFormatter.
# Duplicate A1/A3/A5 into A2/A4/A6 and apply a "0.0" one-decimal number
# format to the duplicated cells so the rounding behaviour can be
# demonstrated/tested.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the values of the odd rows into the even rows directly below them.
$ws.Range("A2").Value2 = 2.04
$ws.Range("A4").Value2 = 2.0499999999999998
$ws.Range("A6").Value2 = 2.06

# Give the mirrored cells a custom "0.0" number format (one decimal place).
$ws.Range("A2").NumberFormat = "0.0"
$ws.Range("A4").NumberFormat = "0.0"
$ws.Range("A6").NumberFormat = "0.0"

# Move the active selection to A6 (was A7).
$ws.Range("A6").Select()

# Page setup as saved by the editing session (A4 paper, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
